$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1940
$ws.Cells.Item(3, 6).Value = 36
$ws.Cells.Item(4, 6).Value = 100
$ws.Cells.Item(5, 6).Value = 423
$ws.Cells.Item(6, 6).Value = 1839
$ws.Cells.Item(7, 6).Value = 860
$ws.Cells.Item(8, 6).Value = 1303
$ws.Cells.Item(9, 6).Value = 542
$ws.Cells.Item(11, 6).Value = 2758
$ws.Cells.Item(13, 6).Value = 868
$ws.Cells.Item(14, 6).Value = 1087
$ws.Cells.Item(15, 6).Value = 580
$ws.Cells.Item(17, 6).Value = 61
$ws.Cells.Item(18, 6).Value = 1584
$ws.Cells.Item(19, 6).Value = 31
$ws.Cells.Item(20, 6).Value = 1237
$ws.Cells.Item(21, 6).Value = 179
$ws.Cells.Item(21, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/i6vAgX8I1719311206769.jpeg"
$ws.Cells.Item(22, 6).Value = 595
$ws.Cells.Item(25, 6).Value = 1447
$ws.Cells.Item(26, 6).Value = 1449
$ws.Cells.Item(27, 6).Value = 1320
$ws.Cells.Item(28, 6).Value = 246
$ws.Cells.Item(29, 6).Value = 1273
$ws.Cells.Item(30, 6).Value = 430
$ws.Cells.Item(32, 6).Value = 957
$ws.Cells.Item(33, 6).Value = 22
$ws.Cells.Item(35, 6).Value = 467
$ws.Cells.Item(36, 6).Value = 36
$ws.Cells.Item(39, 6).Value = 2261
$ws.Cells.Item(40, 6).Value = 138
$ws.Cells.Item(41, 6).Value = 884
$ws.Cells.Item(42, 6).Value = 2754
$ws.Cells.Item(44, 6).Value = 843
$ws.Cells.Item(45, 6).Value = 17

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 60
$ws.Cells.Item(10, 6).Value = 31
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(20, 6).Value = 286
$ws.Cells.Item(22, 6).Value = 275
$ws.Cells.Item(30, 6).Value = 52
$ws.Cells.Item(30, 7).Value = 180
$ws.Cells.Item(31, 6).Value = 216
$ws.Cells.Item(32, 6).Value = 10
$ws.Cells.Item(40, 6).Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 3019
$ws.Cells.Item(6, 6).Value = 4847
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(7, 6).Value = 175
$ws.Cells.Item(9, 6).Value = 658
$ws.Cells.Item(10, 6).Value = 928
$ws.Cells.Item(11, 6).Value = 535
$ws.Cells.Item(12, 6).Value = 606
$ws.Cells.Item(13, 6).Value = 1342
$ws.Cells.Item(14, 6).Value = 378
$ws.Cells.Item(15, 6).Value = 1206

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1940
$ws.Cells.Item(4, 6).Value = 36
$ws.Cells.Item(5, 6).Value = 658
$ws.Cells.Item(6, 6).Value = 928
$ws.Cells.Item(7, 6).Value = 535
$ws.Cells.Item(8, 6).Value = 100
$ws.Cells.Item(9, 6).Value = 606
$ws.Cells.Item(10, 6).Value = 1342
$ws.Cells.Item(11, 6).Value = 423
$ws.Cells.Item(12, 6).Value = 1839
$ws.Cells.Item(13, 6).Value = 860
$ws.Cells.Item(14, 6).Value = 1303
$ws.Cells.Item(16, 6).Value = 542
$ws.Cells.Item(17, 6).Value = 1206
$ws.Cells.Item(18, 6).Value = 2758
$ws.Cells.Item(19, 6).Value = 31
$ws.Cells.Item(21, 6).Value = 868
$ws.Cells.Item(22, 6).Value = 1088
$ws.Cells.Item(23, 6).Value = 580
$ws.Cells.Item(25, 6).Value = 1584
$ws.Cells.Item(26, 6).Value = 31
$ws.Cells.Item(28, 6).Value = 1237
$ws.Cells.Item(29, 6).Value = 179
$ws.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/i6vAgX8I1719311206769.jpeg"
$ws.Cells.Item(30, 6).Value = 595
$ws.Cells.Item(31, 6).Value = 1447
$ws.Cells.Item(32, 6).Value = 1449
$ws.Cells.Item(33, 6).Value = 1320
$ws.Cells.Item(36, 6).Value = 1273
$ws.Cells.Item(37, 6).Value = 430
$ws.Cells.Item(38, 6).Value = 957
$ws.Cells.Item(42, 6).Value = 216
$ws.Cells.Item(43, 6).Value = 2261
$ws.Cells.Item(44, 6).Value = 138
$ws.Cells.Item(45, 6).Value = 884
$ws.Cells.Item(46, 6).Value = 2754
$ws.Cells.Item(47, 6).Value = 843
$ws.Cells.Item(51, 6).Value = 6
